$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2125
$ws.Range("J70").Value = 1250
$ws.Range("L70").Value = 3750
$ws.Range("N70").Value = -4290
$ws.Range("H73").Value = 2125
$ws.Range("J73").Value = 1250
$ws.Range("L73").Value = 3750
$ws.Range("N73").Value = -5622
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = ""
$ws.Range("H113").Value = 2873.5
$ws.Range("I113").Value = 2187.25
$ws.Range("J113").Value = 4246
$ws.Range("K113").Value = 2187.25
$ws.Range("L113").Value = 4246
$ws.Range("M113").Value = 1066.75
$ws.Range("N113").Value = -10754
$ws.Range("H116").Value = 8000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = -14884
$ws.Range("H121").Value = 4420
$ws.Range("J121").Value = 4420
$ws.Range("L121").Value = 13260
$ws.Range("N121").Value = -16754
$ws.Range("H132").Value = 3853.4
$ws.Range("I132").Value = 3853.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11560.2
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9030.200000000001
$ws.Range("N132").Value = ""
$ws.Range("H138").Value = 2097.5557
$ws.Range("I138").Value = 1621.1875
$ws.Range("J138").Value = 5908.5
$ws.Range("K138").Value = 4863.5625
$ws.Range("L138").Value = 17725.5
$ws.Range("M138").Value = 276.4375
$ws.Range("N138").Value = -28005.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1918.3334
$ws.Range("I74").Value = 1702
$ws.Range("K74").Value = 1702
$ws.Range("M74").Value = -828
$ws.Range("H77").Value = 1918.3334
$ws.Range("I77").Value = 1702
$ws.Range("K77").Value = 8510
$ws.Range("M77").Value = -4142
$ws.Range("H122").Value = 1729.5385
$ws.Range("I122").Value = 1346.2222
$ws.Range("J122").Value = 2592
$ws.Range("K122").Value = 4038.6666
$ws.Range("L122").Value = 7776
$ws.Range("M122").Value = -1588.6666
$ws.Range("N122").Value = -12676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1241.1428
$ws.Range("I86").Value = 1091.8572
$ws.Range("K86").Value = 1091.8572
$ws.Range("M86").Value = 31.14280000000008
$ws.Range("H89").Value = 1241.1428
$ws.Range("I89").Value = 1091.8572
$ws.Range("K89").Value = 5459.286
$ws.Range("M89").Value = 156.7139999999999
$ws.Range("H102").Value = 49333
$ws.Range("I102").Value = 49333
$ws.Range("K102").Value = 49333
$ws.Range("M102").Value = -46088
$ws.Range("H134").Value = 2616
$ws.Range("I134").Value = 2616
$ws.Range("K134").Value = 7848
$ws.Range("M134").Value = -5313

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2108
$ws.Range("I58").Value = 1495.5
$ws.Range("K58").Value = 1495.5
$ws.Range("M58").Value = -1292.5
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").Value = ""
$ws.Range("H134").Value = 2949.75
$ws.Range("I134").Value = 2949.75
$ws.Range("K134").Value = 8849.25
$ws.Range("M134").Value = -6314.25
$ws.Range("H136").Value = 2108
$ws.Range("I136").Value = 1495.5
$ws.Range("K136").Value = 4486.5
$ws.Range("M136").Value = -1936.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = ""
$ws.Range("H26").Value = 1633
$ws.Range("J26").Value = 1999
$ws.Range("L26").Value = 5997
$ws.Range("N26").Value = -6573
$ws.Range("H29").Value = 694
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 694
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 2082
$ws.Range("M29").Value = ""
$ws.Range("N29").Value = -2636
$ws.Range("H40").Value = 532.125
$ws.Range("I40").Value = 74.25
$ws.Range("K40").Value = 297
$ws.Range("M40").Value = -228
$ws.Range("H81").Value = 2111.6667
$ws.Range("J81").Value = 2111.6667
$ws.Range("L81").Value = 6335.000100000001
$ws.Range("N81").Value = -8581.000100000001
$ws.Range("H84").Value = 2111.6667
$ws.Range("J84").Value = 2111.6667
$ws.Range("L84").Value = 19005.0003
$ws.Range("N84").Value = -30237.0003
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("H104").Value = 6997
$ws.Range("J104").Value = 6997
$ws.Range("L104").Value = 20991
$ws.Range("N104").Value = -26233
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 4000
$ws.Range("K70").Value = 4000
$ws.Range("M70").Value = -3730
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 4000
$ws.Range("K73").Value = 4000
$ws.Range("M73").Value = -3064
$ws.Range("H99").Value = 9658
$ws.Range("I99").Value = 6599.3335
$ws.Range("J99").Value = 11952
$ws.Range("K99").Value = 6599.3335
$ws.Range("L99").Value = 11952
$ws.Range("M99").Value = -4353.3335
$ws.Range("N99").Value = -16444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2067.5
$ws.Range("I16").Value = 2067.5
$ws.Range("K16").Value = 2067.5
$ws.Range("M16").Value = -1897.5
$ws.Range("H22").Value = 3883.9167
$ws.Range("J22").Value = 3564
$ws.Range("L22").Value = 3564
$ws.Range("N22").Value = -4154
$ws.Range("H27").Value = 3883.9167
$ws.Range("J27").Value = 3564
$ws.Range("L27").Value = 3564
$ws.Range("N27").Value = -3778
$ws.Range("H40").Value = 4524
$ws.Range("I40").Value = 4265.6665
$ws.Range("K40").Value = 4265.6665
$ws.Range("M40").Value = -4129.6665
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = ""
$ws.Range("H93").Value = 1779
$ws.Range("I93").Value = 1776
$ws.Range("K93").Value = 1776
$ws.Range("M93").Value = -528
$ws.Range("H99").Value = 64250
$ws.Range("I99").Value = 64250
$ws.Range("K99").Value = 64250
$ws.Range("M99").Value = -61255
$ws.Range("H100").Value = 3573.9443
$ws.Range("I100").Value = 3729
$ws.Range("J100").Value = 2798.6667
$ws.Range("K100").Value = 3729
$ws.Range("L100").Value = 2798.6667
$ws.Range("M100").Value = -3188
$ws.Range("N100").Value = -3880.6667
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = ""
$ws.Range("H122").Value = 7134.654
$ws.Range("I122").Value = 6656.4
$ws.Range("K122").Value = 19969.2
$ws.Range("M122").Value = -17519.2
$ws.Range("H132").Value = 3546.1
$ws.Range("I132").Value = 3329.3333
$ws.Range("K132").Value = 9987.999899999999
$ws.Range("M132").Value = -7457.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17555.6
$ws.Range("J41").Value = 16950
$ws.Range("L41").Value = 16950
$ws.Range("N41").Value = -17730
$ws.Range("H132").Value = 1400.75
$ws.Range("I132").Value = 1400.75
$ws.Range("K132").Value = 4202.25
$ws.Range("M132").Value = -1672.25
$ws.Range("H136").Value = 729.5714
$ws.Range("I136").Value = 687.36365
$ws.Range("K136").Value = 2062.09095
$ws.Range("M136").Value = 487.9090500000002

Write-Output "Applied all cell updates"